$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1189.25
$ws.Cells.Item(2, 10).Value = 988.6667
$ws.Cells.Item(2, 12).Value = 988.6667
$ws.Cells.Item(2, 14).Value = -1214.6667
$ws.Cells.Item(32, 8).Value = 471078.2
$ws.Cells.Item(32, 9).Value = 527220.4399999999
$ws.Cells.Item(32, 10).Value = 21940
$ws.Cells.Item(32, 11).Value = 527220.4399999999
$ws.Cells.Item(32, 12).Value = 21940
$ws.Cells.Item(32, 13).Value = -526933.4399999999
$ws.Cells.Item(32, 14).Value = -22514
$ws.Cells.Item(61, 8).Value = 7409689.5
$ws.Cells.Item(61, 9).Value = 16667859
$ws.Cells.Item(61, 10).Value = 3154.04
$ws.Cells.Item(61, 11).Value = 16667859
$ws.Cells.Item(61, 12).Value = 3154.04
$ws.Cells.Item(61, 13).Value = -16667647
$ws.Cells.Item(61, 14).Value = -3578.04
$ws.Cells.Item(74, 8).Value = 841.90247
$ws.Cells.Item(74, 9).Value = 556.86664
$ws.Cells.Item(74, 10).Value = 1006.3461
$ws.Cells.Item(74, 11).Value = 556.86664
$ws.Cells.Item(74, 12).Value = 1006.3461
$ws.Cells.Item(74, 13).Value = 317.13336
$ws.Cells.Item(74, 14).Value = -2754.3461
$ws.Cells.Item(77, 8).Value = 841.90247
$ws.Cells.Item(77, 9).Value = 556.86664
$ws.Cells.Item(77, 10).Value = 1006.3461
$ws.Cells.Item(77, 11).Value = 2784.3332
$ws.Cells.Item(77, 12).Value = 5031.7305
$ws.Cells.Item(77, 13).Value = 1583.6668
$ws.Cells.Item(77, 14).Value = -13767.7305
$ws.Cells.Item(116, 8).Value = 1189.25
$ws.Cells.Item(116, 10).Value = 988.6667
$ws.Cells.Item(116, 12).Value = 988.6667
$ws.Cells.Item(116, 14).Value = -5576.6667
$ws.Cells.Item(136, 8).Value = 7409689.5
$ws.Cells.Item(136, 9).Value = 16667859
$ws.Cells.Item(136, 10).Value = 3154.04
$ws.Cells.Item(136, 11).Value = 50003577
$ws.Cells.Item(136, 12).Value = 9462.119999999999
$ws.Cells.Item(136, 13).Value = -50001027
$ws.Cells.Item(136, 14).Value = -14562.12

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1189.25
$ws.Cells.Item(3, 10).Value = 988.6667
$ws.Cells.Item(3, 12).Value = 988.6667
$ws.Cells.Item(3, 14).Value = -1216.6667
$ws.Cells.Item(60, 8).Value = 20499.5
$ws.Cells.Item(60, 10).Value = 20499.5
$ws.Cells.Item(60, 12).Value = 20499.5
$ws.Cells.Item(60, 14).Value = -21697.5
$ws.Cells.Item(134, 8).Value = 3232.5925
$ws.Cells.Item(134, 9).Value = 4020.182
$ws.Cells.Item(134, 10).Value = 2691.125
$ws.Cells.Item(134, 11).Value = 12060.546
$ws.Cells.Item(134, 12).Value = 8073.375
$ws.Cells.Item(134, 13).Value = -9525.545999999998
$ws.Cells.Item(134, 14).Value = -13143.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4153.95
$ws.Cells.Item(31, 9).Value = 1264.9412
$ws.Cells.Item(31, 11).Value = 1264.9412
$ws.Cells.Item(31, 13).Value = -969.9412
$ws.Cells.Item(34, 8).Value = 4153.95
$ws.Cells.Item(34, 9).Value = 1264.9412
$ws.Cells.Item(34, 11).Value = 1264.9412
$ws.Cells.Item(34, 13).Value = -1062.9412
$ws.Cells.Item(58, 8).Value = 2303.2
$ws.Cells.Item(58, 10).Value = 576
$ws.Cells.Item(58, 12).Value = 576
$ws.Cells.Item(58, 14).Value = -982
$ws.Cells.Item(86, 8).Value = 1529.7046
$ws.Cells.Item(86, 9).Value = 1584
$ws.Cells.Item(86, 10).Value = 1413.3572
$ws.Cells.Item(86, 11).Value = 1584
$ws.Cells.Item(86, 12).Value = 1413.3572
$ws.Cells.Item(86, 13).Value = -461
$ws.Cells.Item(86, 14).Value = -3659.3572
$ws.Cells.Item(89, 8).Value = 1529.7046
$ws.Cells.Item(89, 9).Value = 1584
$ws.Cells.Item(89, 10).Value = 1413.3572
$ws.Cells.Item(89, 11).Value = 7920
$ws.Cells.Item(89, 12).Value = 7066.786
$ws.Cells.Item(89, 13).Value = -2304
$ws.Cells.Item(89, 14).Value = -18298.786
$ws.Cells.Item(136, 8).Value = 2303.2
$ws.Cells.Item(136, 10).Value = 576
$ws.Cells.Item(136, 12).Value = 1728
$ws.Cells.Item(136, 14).Value = -6828

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1082.3768
$ws.Cells.Item(68, 9).Value = 1017.2727
$ws.Cells.Item(68, 10).Value = 1094.7241
$ws.Cells.Item(68, 11).Value = 3051.8181
$ws.Cells.Item(68, 12).Value = 3284.1723
$ws.Cells.Item(68, 13).Value = -2240.8181
$ws.Cells.Item(68, 14).Value = -4906.1723
$ws.Cells.Item(71, 8).Value = 1082.3768
$ws.Cells.Item(71, 9).Value = 1017.2727
$ws.Cells.Item(71, 10).Value = 1094.7241
$ws.Cells.Item(71, 11).Value = 9155.454299999999
$ws.Cells.Item(71, 12).Value = 9852.516899999999
$ws.Cells.Item(71, 13).Value = -5099.454299999999
$ws.Cells.Item(71, 14).Value = -17964.5169
$ws.Cells.Item(112, 8).Value = 6200
$ws.Cells.Item(112, 10).Value = 6666.6665
$ws.Cells.Item(112, 12).Value = 19999.9995
$ws.Cells.Item(112, 14).Value = -22215.9995
$ws.Cells.Item(121, 8).Value = 1024.4166
$ws.Cells.Item(121, 9).Value = 409.75
$ws.Cells.Item(121, 10).Value = 1229.3055
$ws.Cells.Item(121, 11).Value = 1229.25
$ws.Cells.Item(121, 12).Value = 3687.9165
$ws.Cells.Item(121, 13).Value = 80.75
$ws.Cells.Item(121, 14).Value = -6307.916499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1695.6471
$ws.Cells.Item(122, 9).Value = 1355.4546
$ws.Cells.Item(122, 10).Value = 2319.3333
$ws.Cells.Item(122, 11).Value = 4066.3638
$ws.Cells.Item(122, 12).Value = 6957.999899999999
$ws.Cells.Item(122, 13).Value = -1616.3638
$ws.Cells.Item(122, 14).Value = -11857.9999
$ws.Cells.Item(132, 8).Value = 1783.7
$ws.Cells.Item(132, 9).Value = 834
$ws.Cells.Item(132, 10).Value = 3999.6667
$ws.Cells.Item(132, 11).Value = 2502
$ws.Cells.Item(132, 12).Value = 11999.0001
$ws.Cells.Item(132, 13).Value = 28
$ws.Cells.Item(132, 14).Value = -17059.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(70, 8).Value = 95000
$ws.Cells.Item(70, 10).Value = 95000
$ws.Cells.Item(70, 12).Value = 95000
$ws.Cells.Item(70, 14).Value = -95540
$ws.Cells.Item(73, 8).Value = 95000
$ws.Cells.Item(73, 10).Value = 95000
$ws.Cells.Item(73, 12).Value = 95000
$ws.Cells.Item(73, 14).Value = -96872
$ws.Cells.Item(132, 8).Value = 3310.0908
$ws.Cells.Item(132, 9).Value = 3354.6365
$ws.Cells.Item(132, 10).Value = 3265.5454
$ws.Cells.Item(132, 11).Value = 10063.9095
$ws.Cells.Item(132, 12).Value = 9796.636200000001
$ws.Cells.Item(132, 13).Value = -7533.9095
$ws.Cells.Item(132, 14).Value = -14856.6362
$ws.Cells.Item(136, 8).Value = 1570.5652
$ws.Cells.Item(136, 9).Value = 1607.3572
$ws.Cells.Item(136, 10).Value = 1513.3334
$ws.Cells.Item(136, 11).Value = 4822.071599999999
$ws.Cells.Item(136, 12).Value = 4540.0002
$ws.Cells.Item(136, 13).Value = -2272.071599999999
$ws.Cells.Item(136, 14).Value = -9640.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(63, 8).Value = 75832.664
$ws.Cells.Item(63, 10).Value = 75832.664
$ws.Cells.Item(63, 12).Value = 75832.664
$ws.Cells.Item(63, 14).Value = -77080.664
$ws.Cells.Item(66, 8).Value = 75832.664
$ws.Cells.Item(66, 10).Value = 75832.664
$ws.Cells.Item(66, 12).Value = 227497.992
$ws.Cells.Item(66, 14).Value = -233737.992
$ws.Cells.Item(81, 8).Value = 4789.0835
$ws.Cells.Item(81, 9).Value = 5928.3335
$ws.Cells.Item(81, 10).Value = 3649.8333
$ws.Cells.Item(81, 11).Value = 11856.667
$ws.Cells.Item(81, 12).Value = 7299.6666
$ws.Cells.Item(81, 13).Value = -10795.667
$ws.Cells.Item(81, 14).Value = -9421.6666
$ws.Cells.Item(84, 8).Value = 4789.0835
$ws.Cells.Item(84, 9).Value = 5928.3335
$ws.Cells.Item(84, 10).Value = 3649.8333
$ws.Cells.Item(84, 11).Value = 59283.335
$ws.Cells.Item(84, 12).Value = 36498.333
$ws.Cells.Item(84, 13).Value = -53979.335
$ws.Cells.Item(84, 14).Value = -47106.333
$ws.Cells.Item(136, 8).Value = 1862.6533
$ws.Cells.Item(136, 9).Value = 1584.5
$ws.Cells.Item(136, 11).Value = 4753.5
$ws.Cells.Item(136, 13).Value = -2203.5
